$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '41.456.12'
$ws.Range('E2').Value = '  -3.31%  '

Set-TextValue $ws.Range('D3') '2.468.48'
$ws.Range('E3').Value = '  -2.60%  '

Set-TextValue $ws.Range('D5') '311.86'
$ws.Range('E5').Value = '  -0.15%  '

Set-TextValue $ws.Range('D6') '94.83'
$ws.Range('E6').Value = '  -5.82%  '

$ws.Range('E7').Value = '  -2.64%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('E9').Value = '  -4.46%  '

Set-TextValue $ws.Range('D10') '33.55'
$ws.Range('E10').Value = '  -6.33%  '

Set-TextValue $ws.Range('D11') '0.0781'
$ws.Range('E11').Value = '  -3.21%  '

$ws.Range('E12').Value = '  -1.04%  '

Set-TextValue $ws.Range('D13') '7.02'
$ws.Range('E13').Value = '  -4.35%  '

Set-TextValue $ws.Range('D14') '2.848.60'
$ws.Range('E14').Value = '  -2.58%  '

Set-TextValue $ws.Range('D15') '2.476.20'
$ws.Range('E15').Value = '  -1.83%  '

$ws.Range('E16').Value = '  -2.98%  '

$ws.Range('E17').Value = '  -3.71%  '

Set-TextValue $ws.Range('D18') '41.414.88'
$ws.Range('E18').Value = '  -3.36%  '

Set-TextValue $ws.Range('D19') '6.33'
$ws.Range('E19').Value = '  -5.27%  '

Set-TextValue $ws.Range('D20') '0.0₃0922'
$ws.Range('E20').Value = '  -3.45%  '

Set-TextValue $ws.Range('D21') '11.24'
$ws.Range('E21').Value = '  -9.13%  '

Set-TextValue $ws.Range('D22') '68.57'
$ws.Range('E22').Value = '  -1.86%  '

Set-TextValue $ws.Range('D23') '236.76'
$ws.Range('E23').Value = '  -3.03%  '

Set-TextValue $ws.Range('D24') '2.75'
$ws.Range('E24').Value = '  -4.94%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D25') '1.00'
$ws.Range('E25').Value = '  -0.04%  '

$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D26') '1.90'
$ws.Range('E26').Value = '  -6.80%  '

Set-TextValue $ws.Range('D27') '24.14'
$ws.Range('E27').Value = '  -5.69%  '

Set-TextValue $ws.Range('D28') '2.21'
$ws.Range('E28').Value = '  -5.40%  '

Set-TextValue $ws.Range('D29') '9.64'
$ws.Range('E29').Value = '  -5.83%  '

Set-TextValue $ws.Range('D30') '36.61'
$ws.Range('E30').Value = '  -6.01%  '

Set-TextValue $ws.Range('D31') '152.05'
$ws.Range('E31').Value = '  -5.70%  '

Set-TextValue $ws.Range('D32') '5.49'
$ws.Range('E32').Value = '  -5.97%  '

Set-TextValue $ws.Range('D33') '2.65'
$ws.Range('E33').Value = '  -4.66%  '

Set-TextValue $ws.Range('D34') '2.59'
$ws.Range('E34').Value = '  -2.45%  '

Set-TextValue $ws.Range('D35') '0.0750'
$ws.Range('E35').Value = '  -5.29%  '

$ws.Range('E36').Value = '  -2.76%  '

Set-TextValue $ws.Range('D37') '17.11'
$ws.Range('E37').Value = '  -7.17%  '

Set-TextValue $ws.Range('D38') '1.88'
$ws.Range('E38').Value = '  -4.73%  '

$ws.Range('E39').Value = '  -2.93%  '

Set-TextValue $ws.Range('D40') '0.102'
$ws.Range('E40').Value = '  -8.14%  '

Set-TextValue $ws.Range('D41') '4.26'
$ws.Range('E41').Value = '  +1.85%  '

$ws.Range('E42').Value = '  +0.16%  '

Set-TextValue $ws.Range('D43') '19.91'
$ws.Range('E43').Value = '  -9.47%  '

Set-TextValue $ws.Range('D44') '1.988.19'
$ws.Range('E44').Value = '  -0.74%  '

Set-TextValue $ws.Range('D45') '0.0286'
$ws.Range('E45').Value = '  -4.54%  '

Set-TextValue $ws.Range('D46') '3.04'
$ws.Range('E46').Value = '  -9.11%  '

Set-TextValue $ws.Range('D47') '8.73'
$ws.Range('E47').Value = '  -5.47%  '

Set-TextValue $ws.Range('D48') '2.714.44'
$ws.Range('E48').Value = '  -2.24%  '

Set-TextValue $ws.Range('D49') '69.90'
$ws.Range('E49').Value = '  -3.66%  '

Set-TextValue $ws.Range('D50') '96.53'
$ws.Range('E50').Value = '  -4.94%  '

Set-TextValue $ws.Range('D51') '74.84'
$ws.Range('E51').Value = '  -6.11%  '
